$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest cryptos data.
# D-column values are forced to remain plain text (matching the original
# inline-string cell type) by using a leading quote, then the cell style
# is reset to Normal so no stray number-format/quote-prefix style sticks.

$ws.Range("D2").Value = "'26.134.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "'1.667.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").Value = "'209.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").Value = "'0.5212"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").Value = "'0.2621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.05%  "

$ws.Range("D9").Value = "'0.06339"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").Value = "'21.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").Value = "'0.07526"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").Value = "'1.667.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "'4.433"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").Value = "'0.5491"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.89%  "

$ws.Range("D15").Value = "'66.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").Value = "'0.000007960"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.56%  "

$ws.Range("D17").Value = "'26.142.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "'1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").Value = "'4.735"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "

$ws.Range("D20").Value = "'186.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").Value = "'10.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.18%  "

$ws.Range("D22").Value = "'6.179"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").Value = "'149.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'0.1248"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").Value = "'7.488"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.01%  "

$ws.Range("D27").Value = "'15.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").Value = "'0.06392"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.57%  "

$ws.Range("D31").Value = "'3.498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.36%  "

$ws.Range("D32").Value = "'3.408"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.33%  "

$ws.Range("D33").Value = "'1.643"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("D34").Value = "'1.004"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.95%  "

$ws.Range("D35").Value = "'2.406"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

$ws.Range("D36").Value = "'0.6010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.26%  "

$ws.Range("D37").Value = "'2.744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").Value = "'1.110.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").Value = "'6.122"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("D40").Value = "'0.01614"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("D41").Value = "'0.8663"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("E42").Value = "  -0.83%  "

$ws.Range("D43").Value = "'100.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "'1.819.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "'55.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.96%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").Value = "'8.049"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("E49").Value = "  -0.95%  "

$ws.Range("D50").Value = "'0.4243"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "'5.920"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
